# Generate Report for Handoff
#
# Updates the localization-status report to reflect that the handoff
# package is now ready: the per-language "Status" cells move from
# "In Translation" to "Ready for handoff", and the associated timestamp
# cells are bumped to the new generation time.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-language status + latest HO xliff generate date ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-18 02:33:38"

# --- zh-cn sheet: status + latest handoff datetime ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-18 02:33:34"

# --- de-de sheet: status + latest handoff datetime ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-18 02:33:38"

# Widen the status columns so the longer "Ready for handoff" label fits
# (mirrors Excel's own autofit behaviour when the report is regenerated).
$overview.Columns("E:F").AutoFit() | Out-Null
$zhcn.Columns("C:C").AutoFit() | Out-Null
$dede.Columns("C:C").AutoFit() | Out-Null
